$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.441.38'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '2.053.94'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'230.74"
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'57.31"
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').Value = "'0.387"
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('D10').Value = "'0.0816"
$ws.Range('E10').Value = '  +3.96%  '
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').Value = "'14.77"
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '2.358.77'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = "'20.93"
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('D15').Value = "'0.762"
$ws.Range('E15').Value = '  -2.30%  '
$ws.Range('D16').Value = "'5.33"
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '2.050.71'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '37.322.30'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').Value = "'6.11"
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').Value = "'69.97"
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').Value = '0.0₃0842'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').Value = "'227.12"
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'2.39"
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = "'2.30"
$ws.Range('E25').Value = '  -4.24%  '
$ws.Range('D26').Value = "'9.61"
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('D27').Value = "'168.22"
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = "'1.41"
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = "'0.130"
$ws.Range('E29').Value = '  -4.48%  '
$ws.Range('D30').Value = "'19.02"
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').Value = "'0.119"
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').Value = "'4.57"
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('D33').Value = "'4.63"
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = "'0.0617"
$ws.Range('E34').Value = '  -2.70%  '
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = "'1.82"
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -4.30%  '
$ws.Range('D39').Value = "'5.39"
$ws.Range('E39').Value = '  -1.63%  '
$ws.Range('D40').Value = "'0.0223"
$ws.Range('E40').Value = '  -5.42%  '
$ws.Range('D41').Value = "'17.20"
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('D42').Value = '1.491.54'
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = "'0.0948"
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'96.97"
$ws.Range('E45').Value = '  -5.46%  '
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('D47').Value = "'1.03"
$ws.Range('E47').Value = '  -3.65%  '
$ws.Range('D48').Value = "'7.19"
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = "'2.93"
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').Value = "'3.72"
$ws.Range('E50').Value = '  -9.31%  '
$ws.Range('D51').Value = '2.244.17'
$ws.Range('E51').Value = '  -1.25%  '
